$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9, shifting existing rows 9-67 down to become 10-68
$ws.Rows.Item(9).Insert()

# Copy the (now shifted) row 10 as a template into the new blank row 9
$ws.Range("A10:R10").Copy()
$ws.Range("A9").PasteSpecial()
$excel.CutCopyMode = 0

# Update the new row 9 with the latest observation's data
$ws.Range("D9").Value = 44537
$ws.Range("J9").Value = 760
$ws.Range("K9").Value = 3500
$ws.Range("L9").Value = 4000
$ws.Range("M9").Value = 3750
$ws.Range("P9").Value = 3750
